$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they stay strings
# (matching the original inlineStr cell type) instead of being converted to numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '43.754.49'
$ws.Range("E2").Value = '  +4.80%  '
$ws.Range("D3").Value = '2.272.76'
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.51%  '
$ws.Range("D5").Value = '232.41'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("D7").Value = '61.44'
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.412'
$ws.Range("E9").Value = '  +2.81%  '
$ws.Range("D10").Value = '0.0920'
$ws.Range("E10").Value = '  +3.42%  '
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '2.584.95'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '15.74'
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").Value = '22.66'
$ws.Range("E14").Value = '  +4.49%  '
$ws.Range("D15").Value = '5.71'
$ws.Range("E15").Value = '  +2.70%  '
$ws.Range("D16").Value = '0.811'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("D17").Value = '2.224.22'
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '43.321.56'
$ws.Range("E18").Value = '  +4.32%  '
$ws.Range("D19").Value = '0.0₃0941'
$ws.Range("E19").Value = '  +4.93%  '
$ws.Range("D20").Value = '6.24'
$ws.Range("E20").Value = '  +3.28%  '
$ws.Range("D21").Value = '73.05'
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("D22").Value = '248.61'
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("E23").Value = '  +8.15%  '
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("D26").Value = '9.83'
$ws.Range("E26").Value = '  +2.72%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.143'
$ws.Range("E27").Value = '  +2.08%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '169.69'
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").Value = '1.50'
$ws.Range("E29").Value = '  +6.09%  '
$ws.Range("D30").Value = '20.61'
$ws.Range("E30").Value = '  +3.37%  '
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("D33").Value = '5.05'
$ws.Range("E33").Value = '  +2.43%  '
$ws.Range("D34").Value = '4.75'
$ws.Range("E34").Value = '  +2.93%  '
$ws.Range("D35").Value = '0.0660'
$ws.Range("E35").Value = '  +6.01%  '
$ws.Range("D36").Value = '6.49'
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("D37").Value = '2.41'
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").Value = '3.62'
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").Value = '0.0251'
$ws.Range("E39").Value = '  +5.35%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("B41").Value = 'TerraClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D41").Value = '0.000229'
$ws.Range("E41").Value = '  -6.47%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.70'
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("D43").Value = '0.0977'
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = '1.22'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").Value = '4.43'
$ws.Range("E45").Value = '  -8.90%  '
$ws.Range("D46").Value = '97.71'
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("D47").Value = '1.470.57'
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("D48").Value = '16.71'
$ws.Range("E48").Value = '  +1.65%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '1.09'
$ws.Range("E49").Value = '  +0.85%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '2.76'
$ws.Range("E50").Value = '  -1.58%  '
$ws.Range("D51").Value = '2.25'
$ws.Range("E51").Value = '  +7.05%  '
